$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 373, shifting existing rows 373-461 down to 374-462.
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new record's data.
$ws.Cells.Item(373, 1).Value2 = 6
$ws.Cells.Item(373, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(373, 3).Value2 = "Metropolitana"
$ws.Cells.Item(373, 4).Value2 = 45204
$ws.Cells.Item(373, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(373, 5).Value2 = 13
$ws.Cells.Item(373, 6).Value2 = 100112026
$ws.Cells.Item(373, 7).Value2 = "Haba"
$ws.Cells.Item(373, 8).Value2 = "Sin especificar"
$ws.Cells.Item(373, 9).Value2 = "Primera"
$ws.Cells.Item(373, 10).Value2 = 600
$ws.Cells.Item(373, 11).Value2 = 10000
$ws.Cells.Item(373, 12).Value2 = 12000
$ws.Cells.Item(373, 13).Value2 = 11167
$ws.Cells.Item(373, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(373, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(373, 16).Value2 = 447
$ws.Cells.Item(373, 17).Value2 = 25
$ws.Cells.Item(373, 18).Value2 = "Hortaliza"
